# Applies the changes described by the diff:
#  - shared string "Y3" -> "Y4" (used by cell N11 "Status" column)
#  - numeric updates to row 11 across several columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Status" text value in N11 from "Y3" to "Y4"
$ws.Range("N11").Value = "Y4"

# Update the numeric metrics in row 11
$ws.Range("A11").Value = 393
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = 391
$ws.Range("D11").Value = 0.5
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 2
$ws.Range("I11").Value = 45
$ws.Range("J11").Value = 177
$ws.Range("M11").Value = 45.3
$ws.Range("Q11").Value = 2
$ws.Range("R11").Value = 0.5
$ws.Range("T11").Value = 45
$ws.Range("U11").Value = 45.3
